# Ajout de la presence de pokemon dans les herbes
# - rename Feuil3 -> pokemon
# - element_terrain (sheet2): add a 4th row giving the % chance to find a
#   pokemon for each terrain element (HERBE=80, SOL=0, ROCHER=0)
# - pokemon (sheet3): new small table describing the two pokemon that can
#   be found (name / image path / rarity %)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # carte
$ws2 = $wb.Worksheets.Item(2)   # element_terrain
$ws3 = $wb.Worksheets.Item(3)   # Feuil3 -> pokemon

# --- rename the 3rd sheet ---
$ws3.Name = "pokemon"

# --- element_terrain: new row 4, "% proba trouver pokemon" ---
$ws2.Range("A4").Value = "% proba trouver pokemon (ex: 60)"
$ws2.Range("A4").Interior.Color = 65535
$ws2.Range("B4").Value = 80
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 0

# --- pokemon sheet content ---
$ws3.Columns.Item(1).ColumnWidth = 34.83

$ws3.Range("A1").Value = "Nom pokemon (ex: Pikatchu)"
$ws3.Range("A1").Interior.Color = 65535
$ws3.Range("B1").Value = "Pikachu"
$ws3.Range("C1").Value = "Mewtwo"

$ws3.Range("A2").Value = "chemin image pkmn (ex: "
$ws3.Range("A2").Interior.Color = 65535
$ws3.Range("B2").Value = "path_pikachu"
$ws3.Range("C2").Value = "path_mewtwo"

$ws3.Range("A3").Value = "rareté (%)"
$ws3.Range("A3").Interior.Color = 65535

# 66.7 / 33.3 are stored as text in the source workbook (not numbers), so
# force the cell to Text before typing the value, then clear the number
# format again so no visible formatting sticks around on the cell.
$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = "66.7"
$ws3.Range("B3").ClearFormats()

$ws3.Range("C3").NumberFormat = "@"
$ws3.Range("C3").Value = "33.3"
$ws3.Range("C3").ClearFormats()

$ws3.Range("A4").Interior.Color = 65535
$ws3.Range("A5").Interior.Color = 65535

# --- selections / active sheet to match the saved view state ---
$ws1.Activate()
$ws1.Range("O19").Select() | Out-Null

$ws3.Activate()
$ws3.Range("D11").Select() | Out-Null

$ws2.Activate()
$ws2.Range("B4").Select() | Out-Null

Write-Output "edit applied"
